$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 17:25"
$ws.Range("B4").Value = 190022
$ws.Range("C4").Value = 1492
$ws.Range("D4").Value = 7274
$ws.Range("E4").Value = 178646
$ws.Range("G4").Value = 49
$ws.Range("H4").Value = 4102
$ws.Range("B16").Value = 10553
$ws.Range("C16").Value = 373
$ws.Range("E16").Value = 8971
$ws.Range("B20").Value = 5907
$ws.Range("C20").Value = 190
$ws.Range("E20").Value = 5576
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 204
$ws.Range("A43").Value = "Grecia"
$ws.Range("B43").Value = 1415
$ws.Range("C43").Value = 101
$ws.Range("D43").Value = 52
$ws.Range("E43").Value = 1313
$ws.Range("F43").Value = 90
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 50
$ws.Range("A44").Value = "Sudafrica"
$ws.Range("B44").Value = 1353
$ws.Range("D44").Value = 50
$ws.Range("E44").Value = 1298
$ws.Range("F44").Value = 7
$ws.Range("H44").Value = 5
$ws.Range("D71").Value = 43
$ws.Range("E71").Value = 422
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 14
$ws.Range("A74").Value = "Moldavia"
$ws.Range("B74").Value = 423
$ws.Range("C74").Value = 70
$ws.Range("D74").Value = 22
$ws.Range("E74").Value = 397
$ws.Range("F74").Value = 44
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 4
$ws.Range("A75").Value = "Bulgaria"
$ws.Range("B75").Value = 422
$ws.Range("C75").Value = 23
$ws.Range("D75").Value = 20
$ws.Range("E75").Value = 393
$ws.Range("F75").Value = 18
$ws.Range("H75").Value = 9
$ws.Range("A76").Value = "Eslovaquia"
$ws.Range("B76").Value = 400
$ws.Range("C76").Value = 37
$ws.Range("E76").Value = 396
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 1
$ws.Range("A77").Value = "Tunez"
$ws.Range("B77").Value = 394
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 3
$ws.Range("E77").Value = 381
$ws.Range("F77").Value = 10
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 10
$ws.Range("A78").Value = "Principado de Andorra"
$ws.Range("B78").Value = 390
$ws.Range("C78").Value = 14
$ws.Range("D78").Value = 10
$ws.Range("E78").Value = 366
$ws.Range("F78").Value = 12
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 14
$ws.Range("A79").Value = "Kazajistan"
$ws.Range("B79").Value = 375
$ws.Range("C79").Value = 32
$ws.Range("E79").Value = 346
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 3
$ws.Range("A80").Value = "Azerbaiyan"
$ws.Range("B80").Value = 359
$ws.Range("C80").Value = 61
$ws.Range("D80").Value = 26
$ws.Range("E80").Value = 328
$ws.Range("F80").Value = 7
$ws.Range("H80").Value = 5
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 354
$ws.Range("C81").Value = 25
$ws.Range("D81").Value = 12
$ws.Range("E81").Value = 333
$ws.Range("F81").Value = 4
$ws.Range("H81").Value = 9
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 320
$ws.Range("C85").Value = 58
$ws.Range("D85").Value = 28
$ws.Range("E85").Value = 283
$ws.Range("F85").Value = 3
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 9
$ws.Range("A86").Value = "Kuwait"
$ws.Range("B86").Value = 317
$ws.Range("C86").Value = 28
$ws.Range("D86").Value = 80
$ws.Range("E86").Value = 237
$ws.Range("F86").Value = 13
$ws.Range("H86").Value = 0
$ws.Range("A87").Value = "Jordania"
$ws.Range("B87").Value = 274
$ws.Range("D87").Value = 30
$ws.Range("E87").Value = 239
$ws.Range("F87").Value = 5
$ws.Range("H87").Value = 5
$ws.Range("F123").Value = 4
$ws.Range("A142").Value = "Republica de Yibuti"
$ws.Range("B142").Value = 33
$ws.Range("C142").Value = 3
$ws.Range("E142").Value = 33
$ws.Range("H142").Value = 0
$ws.Range("A143").Value = "Guam"
$ws.Range("F143").Value = 0
$ws.Range("A144").Value = "El Salvador"
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 31
$ws.Range("F144").Value = 4
$ws.Range("H144").Value = 1
$ws.Range("A145").Value = "Bermudas"
$ws.Range("B145").Value = 32
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 10
$ws.Range("E145").Value = 22
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0
$ws.Range("A146").Value = "Mali"
$ws.Range("B146").Value = 31
$ws.Range("C146").Value = 3
$ws.Range("E146").Value = 28
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 3
$ws.Range("A158").Value = "Bahamas"
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 1
$ws.Range("H158").Value = 0
$ws.Range("A160").Value = "Birmania"
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 0
$ws.Range("H160").Value = 1
$ws.Range("A163").Value = "Namibia"
$ws.Range("C163").Value = 3
$ws.Range("A164").Value = "Mongolia"
$ws.Range("C164").Value = 2
$ws.Range("A170").Value = "Surinam"
$ws.Range("C170").Value = 1
$ws.Range("A171").Value = "Seychelles"
$ws.Range("C171").Value = 0
$ws.Range("A172").Value = "Mozambique"
$ws.Range("C172").Value = 2
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 10
$ws.Range("A173").Value = "Libia"
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 9
$ws.Range("H173").Value = 0
$ws.Range("A174").Value = "Siria"
$ws.Range("D174").Value = 0
$ws.Range("H174").Value = 2
$ws.Range("A175").Value = "Groenlandia"
$ws.Range("B175").Value = 10
$ws.Range("D175").Value = 2
$ws.Range("E175").Value = 8
$ws.Range("A177").Value = "Suazilandia"
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 9
$ws.Range("A178").Value = "Benin"
$ws.Range("B178").Value = 9
$ws.Range("D178").Value = 1
$ws.Range("A182").Value = "Antigua y Barbuda"
$ws.Range("A183").Value = "Republica del Chad"
$ws.Range("A184").Value = "Angola"
$ws.Range("A185").Value = "Sudan"
$ws.Range("A187").Value = "Liberia"
$ws.Range("C187").Value = 3
$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("C188").Value = 1
$ws.Range("A189").Value = "Cabo Verde"
$ws.Range("D189").Value = 0
$ws.Range("H189").Value = 1
$ws.Range("A190").Value = "San Bartolome"
$ws.Range("D190").Value = 1
$ws.Range("H190").Value = 0
$ws.Range("A192").Value = "Montserrat"
$ws.Range("A193").Value = "Fiyi"
